$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0.1223775669932365
$ws.Range("B2").Value = 0.958997905254364
$ws.Range("C2").Value = 0.01862435974180698
$ws.Range("D2").Value = 0.997705340385437
$ws.Range("A3").Value = 0.03228550404310226
$ws.Range("B3").Value = 0.9933617115020752
$ws.Range("C3").Value = 0.01066349260509014
$ws.Range("D3").Value = 0.9981642961502075
$ws.Range("A4").Value = 0.01493936963379383
$ws.Range("B4").Value = 0.9972929358482361
$ws.Range("C4").Value = 0.004882055800408125
$ws.Range("D4").Value = 0.9990821480751038
$ws.Range("A5").Value = 0.009831191971898079
$ws.Range("B5").Value = 0.9977930188179016
$ws.Range("C5").Value = 0.004934723954647779
$ws.Range("D5").Value = 0.9990821480751038
$ws.Range("A6").Value = 0.007990546524524689
$ws.Range("B6").Value = 0.9981206059455872
$ws.Range("C6").Value = 0.003679408924654126
$ws.Range("D6").Value = 0.9995410442352295
$ws.Range("A7").Value = 0.008193503133952618
$ws.Range("B7").Value = 0.997913658618927
$ws.Range("C7").Value = 0.004622749518603086
$ws.Range("D7").Value = 0.9993880987167358
$ws.Range("A8").Value = 0.007798902690410614
$ws.Range("B8").Value = 0.9980171322822571
$ws.Range("C8").Value = 0.004469004459679127
$ws.Range("D8").Value = 0.9995410442352295
$ws.Range("A9").Value = 0.007138410117477179
$ws.Range("B9").Value = 0.9982413053512573
$ws.Range("C9").Value = 0.005931383464485407
$ws.Range("D9").Value = 0.9993880987167358
$ws.Range("A10").Value = 0.008381451480090618
$ws.Range("B10").Value = 0.997948169708252
$ws.Range("C10").Value = 0.005812073592096567
$ws.Range("D10").Value = 0.9995410442352295
$ws.Range("A11").Value = 0.007036083843559027
$ws.Range("B11").Value = 0.9981550574302673
$ws.Range("C11").Value = 0.005002932157367468
$ws.Range("D11").Value = 0.9995410442352295
$ws.Range("A12").Value = 0.007978073321282864
$ws.Range("B12").Value = 0.9978446960449219
$ws.Range("C12").Value = 0.005956779699772596
$ws.Range("D12").Value = 0.9995410442352295
$ws.Range("A13").Value = 0.006428467109799385
$ws.Range("B13").Value = 0.9983620047569275
$ws.Range("C13").Value = 0.006881711538881063
$ws.Range("D13").Value = 0.9995410442352295
$ws.Range("A14").Value = 0.007843729108572006
$ws.Range("B14").Value = 0.997948169708252
$ws.Range("C14").Value = 0.004265444818884134
$ws.Range("D14").Value = 0.9996940493583679
$ws.Range("A15").Value = 0.007689123973250389
$ws.Range("B15").Value = 0.997999906539917
$ws.Range("C15").Value = 0.004163762554526329
$ws.Range("D15").Value = 0.9995410442352295
$ws.Range("A16").Value = 0.007852815091609955
$ws.Range("B16").Value = 0.9978274703025818
$ws.Range("C16").Value = 0.004234489053487778
$ws.Range("D16").Value = 0.9996940493583679
$ws.Range("A17").Value = 0.007323769386857748
$ws.Range("B17").Value = 0.9980860948562622
$ws.Range("C17").Value = 0.004082055762410164
$ws.Range("D17").Value = 0.9996940493583679
$ws.Range("A18").Value = 0.00681043928489089
$ws.Range("B18").Value = 0.9982413053512573
$ws.Range("C18").Value = 0.004034325480461121
$ws.Range("D18").Value = 0.9996940493583679
$ws.Range("A19").Value = 0.006828949321061373
$ws.Range("B19").Value = 0.9981895685195923
$ws.Range("C19").Value = 0.005054155830293894
$ws.Range("D19").Value = 0.9996940493583679
$ws.Range("A20").Value = 0.007985355332493782
$ws.Range("B20").Value = 0.997948169708252
$ws.Range("C20").Value = 0.004474431741982698
$ws.Range("D20").Value = 0.9995410442352295
$ws.Range("A21").Value = 0.007849578745663166
$ws.Range("B21").Value = 0.997965395450592
$ws.Range("C21").Value = 0.005258076824247837
$ws.Range("D21").Value = 0.9995410442352295
$ws.Range("A22").Value = 0.006594669539481401
$ws.Range("B22").Value = 0.9983102679252625
$ws.Range("C22").Value = 0.00402216799557209
$ws.Range("D22").Value = 0.9996940493583679
$ws.Range("A23").Value = 0.007160215172916651
$ws.Range("B23").Value = 0.9981033205986023
$ws.Range("C23").Value = 0.003767822636291385
$ws.Range("D23").Value = 0.9996940493583679
$ws.Range("A24").Value = 0.007728687953203917
$ws.Range("B24").Value = 0.9978619813919067
$ws.Range("C24").Value = 0.003988831304013729
$ws.Range("D24").Value = 0.9996940493583679
$ws.Range("A25").Value = 0.006970594637095928
$ws.Range("B25").Value = 0.9981550574302673
$ws.Range("C25").Value = 0.004047821275889874
$ws.Range("D25").Value = 0.9996940493583679
$ws.Range("A26").Value = 0.007586246822029352
$ws.Range("B26").Value = 0.9981033205986023
$ws.Range("C26").Value = 0.004579795524477959
$ws.Range("D26").Value = 0.9996940493583679
$ws.Range("A27").Value = 0.007849248126149178
$ws.Range("B27").Value = 0.9977757334709167
$ws.Range("C27").Value = 0.004660735372453928
$ws.Range("D27").Value = 0.9996940493583679
$ws.Range("A28").Value = 0.007322008721530437
$ws.Range("B28").Value = 0.9981033205986023
$ws.Range("C28").Value = 0.004038609098643064
$ws.Range("D28").Value = 0.9996940493583679
$ws.Range("A29").Value = 0.007212950848042965
$ws.Range("B29").Value = 0.9981723427772522
$ws.Range("C29").Value = 0.003310447325929999
$ws.Range("D29").Value = 0.9995410442352295
$ws.Range("A30").Value = 0.006164188962429762
$ws.Range("B30").Value = 0.9983792304992676
$ws.Range("C30").Value = 0.004757424350827932
$ws.Range("D30").Value = 0.9995410442352295
$ws.Range("A31").Value = 0.008193585090339184
$ws.Range("B31").Value = 0.9977585077285767
$ws.Range("C31").Value = 0.005165703129023314
$ws.Range("D31").Value = 0.9996940493583679
$ws.Range("A32").Value = 0.006845239549875259
$ws.Range("B32").Value = 0.9981895685195923
$ws.Range("C32").Value = 0.005096174776554108
$ws.Range("D32").Value = 0.9996940493583679
$ws.Range("A33").Value = 0.007180617656558752
$ws.Range("B33").Value = 0.9980688691139221
$ws.Range("C33").Value = 0.005100996699184179
$ws.Range("D33").Value = 0.9996940493583679
$ws.Range("A34").Value = 0.00675384933128953
$ws.Range("B34").Value = 0.9981895685195923
$ws.Range("C34").Value = 0.005710211582481861
$ws.Range("D34").Value = 0.9996940493583679
$ws.Range("A35").Value = 0.006188265047967434
$ws.Range("B35").Value = 0.9984309673309326
$ws.Range("C35").Value = 0.005278698168694973
$ws.Range("D35").Value = 0.9995410442352295
$ws.Range("A36").Value = 0.006439041811972857
$ws.Range("B36").Value = 0.9983274936676025
$ws.Range("C36").Value = 0.002918129554018378
$ws.Range("D36").Value = 0.9995410442352295
$ws.Range("A37").Value = 0.007952136918902397
$ws.Range("B37").Value = 0.9978446960449219
$ws.Range("C37").Value = 0.004523593932390213
$ws.Range("D37").Value = 0.9995410442352295
$ws.Range("A38").Value = 0.006706835702061653
$ws.Range("B38").Value = 0.9982240200042725
$ws.Range("C38").Value = 0.004524242598563433
$ws.Range("D38").Value = 0.9996940493583679
$ws.Range("A39").Value = 0.007013918831944466
$ws.Range("B39").Value = 0.9981378316879272
$ws.Range("C39").Value = 0.005669512320309877
$ws.Range("D39").Value = 0.9995410442352295
$ws.Range("A40").Value = 0.006953238509595394
$ws.Range("B40").Value = 0.9981723427772522
$ws.Range("C40").Value = 0.005232291761785746
$ws.Range("D40").Value = 0.9995410442352295
$ws.Range("A41").Value = 0.006810517981648445
$ws.Range("B41").Value = 0.9982240200042725
$ws.Range("C41").Value = 0.007534458767622709
$ws.Range("D41").Value = 0.9995410442352295
$ws.Range("A42").Value = 0.007293462287634611
$ws.Range("B42").Value = 0.9982240200042725
$ws.Range("C42").Value = 0.009473095647990704
$ws.Range("D42").Value = 0.9995410442352295
$ws.Range("A43").Value = 0.006847502663731575
$ws.Range("B43").Value = 0.9981895685195923
$ws.Range("C43").Value = 0.006563213188201189
$ws.Range("D43").Value = 0.9995410442352295
$ws.Range("A44").Value = 0.006670957431197166
$ws.Range("B44").Value = 0.9982585310935974
$ws.Range("C44").Value = 0.005262956488877535
$ws.Range("D44").Value = 0.9996940493583679
$ws.Range("A45").Value = 0.007284670602530241
$ws.Range("B45").Value = 0.9981206059455872
$ws.Range("C45").Value = 0.006647361908107996
$ws.Range("D45").Value = 0.9995410442352295
$ws.Range("A46").Value = 0.007388622500002384
$ws.Range("B46").Value = 0.9981206059455872
$ws.Range("C46").Value = 0.004804467782378197
$ws.Range("D46").Value = 0.9995410442352295
$ws.Range("A47").Value = 0.006299341563135386
$ws.Range("B47").Value = 0.9983447194099426
$ws.Range("C47").Value = 0.005138139706104994
$ws.Range("D47").Value = 0.9995410442352295
$ws.Range("A48").Value = 0.007272697985172272
$ws.Range("B48").Value = 0.998051643371582
$ws.Range("C48").Value = 0.005734541453421116
$ws.Range("D48").Value = 0.9995410442352295
$ws.Range("A49").Value = 0.007220827508717775
$ws.Range("B49").Value = 0.9982067942619324
$ws.Range("C49").Value = 0.004502250347286463
$ws.Range("D49").Value = 0.9995410442352295
$ws.Range("A50").Value = 0.006720075383782387
$ws.Range("B50").Value = 0.9981378316879272
$ws.Range("C50").Value = 0.008276755921542645
$ws.Range("D50").Value = 0.9995410442352295
$ws.Range("A51").Value = 0.007204321678727865
$ws.Range("B51").Value = 0.9980343580245972
$ws.Range("C51").Value = 0.008054028265178204
$ws.Range("D51").Value = 0.9995410442352295
